$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 (n): < 65 anos count 160 -> 159
$t.Cell(2, 2).Range.Text = "159"

# Row 3 (Numero.Medicamentos p-value): 0.155 -> 0.141
$t.Cell(3, 4).Range.Text = "0.141"

# Row 5 (Numero.Medicamentos = 1): 49 (30.6) -> 49 (30.8)
$t.Cell(5, 2).Range.Text = "49 (30.8)"

# Row 6 (Numero.Medicamentos = 2): 89 (55.6) -> 88 (55.3)
$t.Cell(6, 2).Range.Text = "88 (55.3)"

# Row 7 (Numero.Medicamentos = 3): 13 ( 8.1) -> 13 ( 8.2)
$t.Cell(7, 2).Range.Text = "13 ( 8.2)"

# Row 9 (Dabigatrana): 2 ( 1.2) -> 2 ( 1.3)
$t.Cell(9, 2).Range.Text = "2 ( 1.3)"

# Row 10 (Enoxaparina): 151 (94.4) -> 150 (94.3)
$t.Cell(10, 2).Range.Text = "150 (94.3)"

# Row 10 (Enoxaparina p-value): 0.288 -> 0.287
$t.Cell(10, 4).Range.Text = "0.287"

# Row 11 (Rivaroxabana): 30 (18.8) -> 30 (18.9)
$t.Cell(11, 2).Range.Text = "30 (18.9)"

# Row 12 (Warfarina): 87 (56.5) -> 86 (56.2)
$t.Cell(12, 2).Range.Text = "86 (56.2)"

# Row 12 (Warfarina p-value): 0.384 -> 0.382
$t.Cell(12, 4).Range.Text = "0.382"
